$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Populate the 3x5 data block coming from the database query ---
for ($r = 1; $r -le 3; $r++) {
    $ws.Cells.Item($r, 1).Value = "boa tarde"
    $ws.Cells.Item($r, 2).Value = 5
    $ws.Cells.Item($r, 3).Value = "ar da sala"
    $ws.Cells.Item($r, 4).Value = 23
    $ws.Cells.Item($r, 5).Value = $true
}

# Column A retained the old single-cell style ("bom dia") from before the
# edit; reset the whole block back to the workbook's Normal style so every
# cell renders with plain/general formatting (text columns general-aligned,
# numbers/boolean using the column defaults).
$ws.Range("A1:E3").Style = "Normal"

# The old sheet had a taller, explicit row height for row 1 left over from
# the previous single-cell layout; auto-fit rows back to the default height.
$ws.Range("1:3").AutoFit()
